$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6226591760299626
$ws1.Range("C2").Value = 0.5746864310148233
$ws1.Range("D2").Value = 0.9438202247191011
$ws1.Range("E2").Value = 0.7143869596031184
$ws1.Range("F2").Value = 0.8363757052771325
$ws1.Range("G2").Value = 0.921065579531876
$ws1.Range("H2").Value = 0.7968410273674761
$ws1.Range("I2").Value = 504
$ws1.Range("J2").Value = 373
$ws1.Range("K2").Value = 161
$ws1.Range("L2").Value = 30

# ---- Sheet: Classification Report ----
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2
$ws2.Range("B2").Value = 0.8429319371727748
$ws2.Range("C2").Value = 0.301498127340824
$ws2.Range("D2").Value = 0.4441379310344827

# row 3
$ws2.Range("B3").Value = 0.5746864310148233
$ws2.Range("C3").Value = 0.9438202247191011
$ws2.Range("D3").Value = 0.7143869596031184

# row 4
$ws2.Range("B4").Value = 0.6226591760299626
$ws2.Range("C4").Value = 0.6226591760299626
$ws2.Range("D4").Value = 0.6226591760299626
$ws2.Range("E4").Value = 0.6226591760299626

# row 5
$ws2.Range("B5").Value = 0.7088091840937991
$ws2.Range("C5").Value = 0.6226591760299626
$ws2.Range("D5").Value = 0.5792624453188006

# row 6
$ws2.Range("B6").Value = 0.7088091840937991
$ws2.Range("C6").Value = 0.6226591760299626
$ws2.Range("D6").Value = 0.5792624453188006

# ---- Sheet: Confusion Matrix ----
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 161
$ws3.Range("C2").Value = 373
$ws3.Range("B3").Value = 30
$ws3.Range("C3").Value = 504
